$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A49").Value = "kontrola a oponentura 4.iterace ŠIP"
$ws.Range("B49").Value = 1

$ws.Range("A50").Select()
